# Applies the scheduled-runner profit recalculation update described in the
# commit diff: per-row currentAveragePrice / LevePrice / LeveProfit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1164.4706
$ws.Range("I19").Value = 833.6667
$ws.Range("J19").Value = 1235.3572
$ws.Range("K19").Value = 833.6667
$ws.Range("L19").Value = 1235.3572
$ws.Range("M19").Value = -658.6667
$ws.Range("N19").Value = -1585.3572

# Row 43
$ws.Range("H43").Value = 2077.3333
$ws.Range("I43").Value = 2360.4167
$ws.Range("J43").Value = 1322.4445
$ws.Range("K43").Value = 2360.4167
$ws.Range("L43").Value = 1322.4445
$ws.Range("M43").Value = -2291.4167
$ws.Range("N43").Value = -1460.4445

# Row 88
$ws.Range("H88").Value = 3762.077
$ws.Range("J88").Value = 3762.077
$ws.Range("L88").Value = 3762.077
$ws.Range("N88").Value = -4574.077

# Row 91
$ws.Range("H91").Value = 3762.077
$ws.Range("J91").Value = 3762.077
$ws.Range("L91").Value = 3762.077
$ws.Range("N91").Value = -6570.077

# Row 94
$ws.Range("H94").Value = 142858140
$ws.Range("I94").Value = 1175
$ws.Range("K94").Value = 1175
$ws.Range("M94").Value = -724

# Row 116
$ws.Range("H116").Value = 2529.7646
$ws.Range("I116").Value = 2800
$ws.Range("J116").Value = 2417.1667
$ws.Range("K116").Value = 2800
$ws.Range("L116").Value = 2417.1667
$ws.Range("M116").Value = 642
$ws.Range("N116").Value = -9301.1667

# Row 132
$ws.Range("H132").Value = 8626515
$ws.Range("I132").Value = 10006117
$ws.Range("J132").Value = 4001.5
$ws.Range("K132").Value = 30018351
$ws.Range("L132").Value = 12004.5
$ws.Range("M132").Value = -30015821
$ws.Range("N132").Value = -17064.5

# Row 135
$ws.Range("H135").Value = 980.3077
$ws.Range("I135").Value = 818.5238000000001
$ws.Range("K135").Value = 7366.7142
$ws.Range("M135").Value = -4831.7142

# Row 141
$ws.Range("H141").Value = 1806.9333
$ws.Range("I141").Value = 1586
$ws.Range("K141").Value = 4758
$ws.Range("M141").Value = 422

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 47167.227
$ws.Range("I102").Value = 68074.336
$ws.Range("J102").Value = 2366.2856
$ws.Range("K102").Value = 68074.336
$ws.Range("L102").Value = 2366.2856
$ws.Range("M102").Value = -66452.336
$ws.Range("N102").Value = -5610.2856

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 93287.5
$ws.Range("I86").Value = 111399
$ws.Range("J86").Value = 2730
$ws.Range("K86").Value = 111399
$ws.Range("L86").Value = 2730
$ws.Range("M86").Value = -110276
$ws.Range("N86").Value = -4976

# Row 89
$ws.Range("H89").Value = 93287.5
$ws.Range("I89").Value = 111399
$ws.Range("J89").Value = 2730
$ws.Range("K89").Value = 556995
$ws.Range("L89").Value = 13650
$ws.Range("M89").Value = -551379
$ws.Range("N89").Value = -24882

# Row 99
$ws.Range("H99").Value = 1931.3334
$ws.Range("I99").Value = 1926.6666
$ws.Range("J99").Value = 1933.6666
$ws.Range("K99").Value = 1926.6666
$ws.Range("L99").Value = 1933.6666
$ws.Range("M99").Value = -428.6666
$ws.Range("N99").Value = -4929.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 36
$ws.Range("H36").Value = 6774
$ws.Range("I36").Value = 6774
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 6774
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -6386
$ws.Range("N36").ClearContents()

# Row 40
$ws.Range("H40").Value = 6774
$ws.Range("I40").Value = 6774
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6774
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6614
$ws.Range("N40").ClearContents()

# Row 45
$ws.Range("H45").Value = 12500
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 94
$ws.Range("H94").Value = 1216.5
$ws.Range("J94").Value = 1216.5
$ws.Range("L94").Value = 1216.5
$ws.Range("N94").Value = -2118.5

# Row 99
$ws.Range("H99").Value = 22798.6
$ws.Range("J99").Value = 27253.25
$ws.Range("L99").Value = 27253.25
$ws.Range("N99").Value = -30249.25

# Row 126
$ws.Range("H126").Value = 22798.6
$ws.Range("J126").Value = 27253.25
$ws.Range("L126").Value = 81759.75
$ws.Range("N126").Value = -86699.75

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 478278.4
$ws.Range("J37").Value = 478278.4
$ws.Range("L37").Value = 1434835.2
$ws.Range("N37").Value = -1435059.2

# Row 131
$ws.Range("H131").Value = 797.26
$ws.Range("I131").Value = 429.9
$ws.Range("J131").Value = 838.07776
$ws.Range("K131").Value = 1289.7
$ws.Range("L131").Value = 2514.23328
$ws.Range("M131").Value = 3750.3
$ws.Range("N131").Value = -12594.23328

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 69300.74000000001
$ws.Range("I70").Value = 91289.83
$ws.Range("J70").Value = 6082.125
$ws.Range("K70").Value = 91289.83
$ws.Range("L70").Value = 6082.125
$ws.Range("M70").Value = -91019.83
$ws.Range("N70").Value = -6622.125

# Row 73
$ws.Range("H73").Value = 69300.74000000001
$ws.Range("I73").Value = 91289.83
$ws.Range("J73").Value = 6082.125
$ws.Range("K73").Value = 91289.83
$ws.Range("L73").Value = 6082.125
$ws.Range("M73").Value = -90353.83
$ws.Range("N73").Value = -7954.125

# Row 102
$ws.Range("H102").Value = 3284.5
$ws.Range("I102").Value = 3427.3572
$ws.Range("J102").Value = 2951.1667
$ws.Range("K102").Value = 3427.3572
$ws.Range("L102").Value = 2951.1667
$ws.Range("M102").Value = -1805.3572
$ws.Range("N102").Value = -6195.1667

# Row 126
$ws.Range("H126").Value = 2771.8572
$ws.Range("I126").Value = 3898.3333
$ws.Range("J126").Value = 2321.2666
$ws.Range("K126").Value = 11694.9999
$ws.Range("L126").Value = 6963.7998
$ws.Range("M126").Value = -9224.999899999999
$ws.Range("N126").Value = -11903.7998

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 72719.14
$ws.Range("I16").Value = 112495.664
$ws.Range("J16").Value = 1121.4
$ws.Range("K16").Value = 112495.664
$ws.Range("L16").Value = 1121.4
$ws.Range("M16").Value = -112325.664
$ws.Range("N16").Value = -1461.4

# Row 136
$ws.Range("H136").Value = 1640.5
$ws.Range("I136").Value = 1600.1666
$ws.Range("J136").Value = 1701
$ws.Range("K136").Value = 4800.4998
$ws.Range("L136").Value = 5103
$ws.Range("M136").Value = -2250.4998
$ws.Range("N136").Value = -10203

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 10000
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10780

# Row 122
$ws.Range("H122").Value = 1670.8
$ws.Range("I122").Value = 1077
$ws.Range("J122").Value = 2066.6667
$ws.Range("K122").Value = 3231
$ws.Range("L122").Value = 6200.000100000001
$ws.Range("M122").Value = -781
$ws.Range("N122").Value = -11100.0001
